$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Resort / correct the LP phases data (columns A, B, D, F) ---
# Row 2
$ws.Range("A2").Value = 53
$ws.Range("B2").Value = "fi"
$ws.Range("D2").Value = "te"
$ws.Range("F2").Value = "se"

# Row 3
$ws.Range("A3").Value = 38
$ws.Range("D3").Value = "la"

# Row 4
$ws.Range("A4").Value = 17
$ws.Range("B4").Value = "pe"
$ws.Range("D4").Value = "pa"
$ws.Range("F4").Value = "to"

# Row 13
$ws.Range("A13").Value = 16
$ws.Range("B13").Value = "pe"
$ws.Range("D13").Value = "pi"
$ws.Range("F13").Value = "to"

# Row 14
$ws.Range("A14").Value = 50
$ws.Range("D14").Value = "sa"

# Row 15
$ws.Range("A15").Value = 30
$ws.Range("D15").Value = "be"

# Row 16
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "pe"
$ws.Range("D16").Value = "le"
$ws.Range("F16").Value = "to"

# Row 17
$ws.Range("A17").Value = 22
$ws.Range("D17").Value = "sa"

# Row 18
$ws.Range("A18").Value = 44
$ws.Range("B18").Value = "fi"
$ws.Range("F18").Value = "se"

# Row 19
$ws.Range("A19").Value = 37
$ws.Range("B19").Value = "fi"
$ws.Range("D19").Value = "ko"
$ws.Range("F19").Value = "se"

# Row 20
$ws.Range("A20").Value = 12
$ws.Range("B20").Value = "pe"
$ws.Range("F20").Value = "to"

# Row 21
$ws.Range("A21").Value = 7
$ws.Range("D21").Value = "ka"

# Row 22
$ws.Range("A22").Value = 40
$ws.Range("B22").Value = "fi"
$ws.Range("D22").Value = "li"
$ws.Range("F22").Value = "se"

# Row 23
$ws.Range("A23").Value = 20
$ws.Range("D23").Value = "ro"

# Row 24
$ws.Range("A24").Value = 56
$ws.Range("B24").Value = "fi"
$ws.Range("D24").Value = "we"
$ws.Range("F24").Value = "se"

# Row 25
$ws.Range("A25").Value = 8
$ws.Range("B25").Value = "pe"
$ws.Range("D25").Value = "ki"
$ws.Range("F25").Value = "to"

# Row 26
$ws.Range("A26").Value = 28
$ws.Range("D26").Value = "we"

# Row 27
$ws.Range("A27").Value = 33
$ws.Range("B27").Value = "fi"
$ws.Range("D27").Value = "fa"
$ws.Range("F27").Value = "se"

# Row 28
$ws.Range("A28").Value = 25
$ws.Range("B28").Value = "pe"
$ws.Range("D28").Value = "te"
$ws.Range("F28").Value = "to"

# Row 29
$ws.Range("A29").Value = 43
$ws.Range("B29").Value = "fi"
$ws.Range("D29").Value = "ni"
$ws.Range("F29").Value = "se"

# Row 30
$ws.Range("A30").Value = 36
$ws.Range("B30").Value = "fi"
$ws.Range("D30").Value = "ki"
$ws.Range("F30").Value = "se"

# Row 31
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "pe"
$ws.Range("D31").Value = "di"
$ws.Range("F31").Value = "to"

# Row 32
$ws.Range("A32").Value = 48
$ws.Range("D32").Value = "ro"

# Row 33
$ws.Range("A33").Value = 9
$ws.Range("D33").Value = "ko"

# Row 34
$ws.Range("A34").Value = 49
$ws.Range("D34").Value = "ri"

# Row 35
$ws.Range("A35").Value = 26
$ws.Range("D35").Value = "ti"

# Row 36
$ws.Range("A36").Value = 35
$ws.Range("D36").Value = "ka"

# Row 37
$ws.Range("A37").Value = 13
$ws.Range("B37").Value = "pe"
$ws.Range("D37").Value = "mi"
$ws.Range("F37").Value = "to"

# --- View changes: zoom to 100%, and move the active selection to C20 ---
$ws.Activate() | Out-Null
$ws.Range("C20").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
